# Applies the "Atualizado por script em 08-11-2023 08:45" update:
#  - Swaps the match-detail columns (F:V) between several pairs of rows
#    (the home/away designation of certain fixtures was corrected while
#    keeping the Indice/pais/torneio/temporada/data_partida columns A:E
#    untouched).
#  - Appends two new rows (162 and 163) describing newly scraped matches.
#
# NOTE: this runtime only reliably binds POSITIONAL function parameters,
# so all helper functions below avoid named-parameter calls.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($Row1, $Row2) {
    $ColStart = 6   # column F
    $ColEnd = 22    # column V
    for ($c = $ColStart; $c -le $ColEnd; $c++) {
        $cell1 = $ws.Cells.Item($Row1, $c)
        $cell2 = $ws.Cells.Item($Row2, $c)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}

# Row pairs whose F:V contents were swapped
Swap-RowData 22 23
Swap-RowData 37 38
Swap-RowData 76 77
Swap-RowData 96 97
Swap-RowData 107 108
Swap-RowData 129 130

function Set-MatchRow(
    $Row, $Indice, $Pais, $Torneio, $Temporada, $DataPartida,
    $Home, $HomeGols, $Away, $AwayGols,
    $HomeOpenOdds, $HomeOpenDataHora, $HomeCloseOdds, $HomeCloseDataHora,
    $DrawOpenOdds, $DrawOpenDataHora, $DrawCloseOdds, $DrawCloseDataHora,
    $AwayOpenOdds, $AwayOpenDataHora, $AwayCloseOdds, $AwayCloseDataHora,
    $Url
) {
    $PrevRow = $Row - 1

    $ws.Cells.Item($Row, 1).Value2 = $Indice
    $ws.Cells.Item($PrevRow, 1).Copy()
    $ws.Cells.Item($Row, 1).PasteSpecial(-4122) # xlPasteFormats

    $ws.Cells.Item($Row, 2).Value2 = $Pais
    $ws.Cells.Item($Row, 3).Value2 = $Torneio
    $ws.Cells.Item($Row, 4).Value2 = $Temporada

    $ws.Cells.Item($Row, 5).Value2 = $DataPartida
    $ws.Cells.Item($PrevRow, 5).Copy()
    $ws.Cells.Item($Row, 5).PasteSpecial(-4122) # xlPasteFormats

    $ws.Cells.Item($Row, 6).Value2 = $Home
    $ws.Cells.Item($Row, 7).Value2 = $HomeGols
    $ws.Cells.Item($Row, 8).Value2 = $Away
    $ws.Cells.Item($Row, 9).Value2 = $AwayGols
    $ws.Cells.Item($Row, 10).Value2 = $HomeOpenOdds
    $ws.Cells.Item($Row, 11).Value2 = $HomeOpenDataHora
    $ws.Cells.Item($Row, 12).Value2 = $HomeCloseOdds
    $ws.Cells.Item($Row, 13).Value2 = $HomeCloseDataHora
    $ws.Cells.Item($Row, 14).Value2 = $DrawOpenOdds
    $ws.Cells.Item($Row, 15).Value2 = $DrawOpenDataHora
    $ws.Cells.Item($Row, 16).Value2 = $DrawCloseOdds
    $ws.Cells.Item($Row, 17).Value2 = $DrawCloseDataHora
    $ws.Cells.Item($Row, 18).Value2 = $AwayOpenOdds
    $ws.Cells.Item($Row, 19).Value2 = $AwayOpenDataHora
    $ws.Cells.Item($Row, 20).Value2 = $AwayCloseOdds
    $ws.Cells.Item($Row, 21).Value2 = $AwayCloseDataHora
    $ws.Cells.Item($Row, 22).Value2 = $Url
}

# New row 162: Persib Bandung 2 x 2 Arema FC
Set-MatchRow 162 161 "indonesia" "liga-1" "2023-2024" 45238.375 `
    "Persib Bandung" 2 "Arema FC" 2 `
    1.33 "06/11/2023 21:13" 1.22 "08/11/2023 08:55" `
    4.88 "06/11/2023 21:13" 5.42 "08/11/2023 08:55" `
    7.04 "06/11/2023 21:13" 9.9 "08/11/2023 08:55" `
    "https://www.betexplorer.com/football/indonesia/liga-1/persib-bandung-arema-fc/MBllPQtC/"

# New row 163: Persik Kediri 4 x 0 Madura United
Set-MatchRow 163 162 "indonesia" "liga-1" "2023-2024" 45238.375 `
    "Persik Kediri" 4 "Madura United" 0 `
    2.53 "06/11/2023 21:11" 2.51 "08/11/2023 08:47" `
    3.14 "06/11/2023 21:11" 3.42 "08/11/2023 08:47" `
    2.56 "06/11/2023 21:11" 2.42 "08/11/2023 08:45" `
    "https://www.betexplorer.com/football/indonesia/liga-1/persik-kediri-madura-united/rLkpQ6R5/"

Write-Host "Edit applied successfully"
